$wb = $excel.ActiveWorkbook

# NOTE: "Vector_bf" and "Vector_BF" only differ by case, and Worksheets.Item(name)
# lookups are case-insensitive, so those two sheets are addressed by their
# (stable) tab index instead of by name to avoid ambiguity.

function Set-TextValue($range, [string]$text) {
    # Force the cell to stay a text value even when $text parses as a number
    # (mirrors the "format cell as Text" step a human would take in Excel
    # before typing a numeric-looking label into it).
    $range.NumberFormat = "@"
    $range.Value = $text
}

# --- Restricciones_del_follower sheet ---
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2 (J_0_L0_v)
$wsFollower.Range("A2").Value = "-22.36259364181055 + 2.777961943082056y"
Set-TextValue $wsFollower.Range("B2") "22.36259364181055"
Set-TextValue $wsFollower.Range("D2") "0.13"
Set-TextValue $wsFollower.Range("E2") "8.100000000000001"
Set-TextValue $wsFollower.Range("F2") "8.299999999999999"

# Row 3 (J_0_L0_v)
$wsFollower.Range("A3").Value = "0.17256554072491337 - x + 0.7363272620217498y"
Set-TextValue $wsFollower.Range("B3") "-3.1725655407249134"
Set-TextValue $wsFollower.Range("D3") "0.6"
Set-TextValue $wsFollower.Range("E3") "1.6"
Set-TextValue $wsFollower.Range("F3") "2.2"

# Row 4 (J_0_LP_v)
$wsFollower.Range("A4").Value = "25.90361075419989 + x - 3.9756038203975015y"
Set-TextValue $wsFollower.Range("B4") "-37.90361075419989"
Set-TextValue $wsFollower.Range("D4") "0.62"
Set-TextValue $wsFollower.Range("E4") "8.0"
Set-TextValue $wsFollower.Range("F4") "1.5"

# Row 5 (J_Ne_L0_v)
$wsFollower.Range("A5").Value = "-42.731732216883216 + 4x + 2.242451207066238y"
Set-TextValue $wsFollower.Range("B5") "30.451732216883215"
Set-TextValue $wsFollower.Range("D5") "0.32"
Set-TextValue $wsFollower.Range("E5") "9.0"
Set-TextValue $wsFollower.Range("F5") "6.7"

# --- Punto_modificado sheet ---
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto.Range("A2") "6.1"
Set-TextValue $wsPunto.Range("B2") "8.05"

# --- Vector_bf sheet (index 5) ---
$wsBf = $wb.Worksheets.Item(5)
Set-TextValue $wsBf.Range("A2") "-0.055641427428462675"

# --- Vector_BF sheet (index 6) ---
$wsBF = $wb.Worksheets.Item(6)
Set-TextValue $wsBF.Range("A2") "-41.4"
Set-TextValue $wsBF.Range("A3") "-9.056845658615586"

# --- Vector_Alpha sheet ---
$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 2.9878019101987507
